# Apply "new basic stats 2022" update to the governance two-var-stats sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) labels: replace leading "/" and embedded
# "/" separators with "-", and spaces within each segment with "_".
$ws.Range("B1").Value = "Government-Cadw"
$ws.Range("C1").Value = "Government-Local_Authority"
$ws.Range("D1").Value = "Government-National"
$ws.Range("E1").Value = "Government-Other"
$ws.Range("F1").Value = "Independent-English_Heritage"
$ws.Range("G1").Value = "Independent-Historic_Environment_Scotland"
$ws.Range("H1").Value = "Independent-National_Trust"
$ws.Range("I1").Value = "Independent-National_Trust_for_Scotland"
$ws.Range("J1").Value = "Independent-Not_for_profit"
$ws.Range("K1").Value = "Independent-Private"
$ws.Range("L1").Value = "Independent-Unknown"
$ws.Range("M1").Value = "University"
$ws.Range("N1").Value = "Unknown"

# --- Update data values (row 2: England)
$ws.Range("C2").Value = 657
$ws.Range("D2").Value = 56
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 53
$ws.Range("H2").Value = 164
$ws.Range("J2").Value = 1341
$ws.Range("K2").Value = 588
$ws.Range("L2").Value = 165
$ws.Range("M2").Value = 76
$ws.Range("N2").Value = 83

# --- Row 3: Northern Ireland
$ws.Range("C3").Value = 34
$ws.Range("J3").Value = 32
$ws.Range("K3").Value = 18

# --- Row 4: Scotland
$ws.Range("C4").Value = 162
$ws.Range("G4").Value = 21
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 27
$ws.Range("J4").Value = 251
$ws.Range("K4").Value = 77
$ws.Range("L4").Value = 39
$ws.Range("M4").Value = 28
$ws.Range("N4").Value = 7

# --- Row 5: Wales
$ws.Range("C5").Value = 61
$ws.Range("D5").Value = 11
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 56
